$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5411.4165
$ws.Range("I132").Value = 6467.579
$ws.Range("K132").Value = 19402.737
$ws.Range("M132").Value = -16872.737

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 818009.25
$ws.Range("I2").Value = 739.95654
$ws.Range("K2").Value = 739.95654
$ws.Range("M2").Value = -626.95654

$ws.Range("H32").Value = 6581.088
$ws.Range("I32").Value = 4870.1816
$ws.Range("K32").Value = 4870.1816
$ws.Range("M32").Value = -4583.1816

$ws.Range("H61").Value = 1323543.4
$ws.Range("I61").Value = 1462771.9
$ws.Range("J61").Value = 873.25
$ws.Range("K61").Value = 1462771.9
$ws.Range("L61").Value = 873.25
$ws.Range("M61").Value = -1462559.9
$ws.Range("N61").Value = -1297.25

$ws.Range("H74").Value = 11908827
$ws.Range("I74").Value = 15625578
$ws.Range("J74").Value = 15225.3
$ws.Range("K74").Value = 15625578
$ws.Range("L74").Value = 15225.3
$ws.Range("M74").Value = -15624704
$ws.Range("N74").Value = -16973.3

$ws.Range("H77").Value = 11908827
$ws.Range("I77").Value = 15625578
$ws.Range("J77").Value = 15225.3
$ws.Range("K77").Value = 78127890
$ws.Range("L77").Value = 76126.5
$ws.Range("M77").Value = -78123522
$ws.Range("N77").Value = -84862.5

$ws.Range("H110").Value = 1252.5
$ws.Range("I110").Value = 1069.4667
$ws.Range("K110").Value = 1069.4667
$ws.Range("M110").Value = 975.5333000000001

$ws.Range("H116").Value = 818009.25
$ws.Range("I116").Value = 739.95654
$ws.Range("K116").Value = 739.95654
$ws.Range("M116").Value = 1554.04346

$ws.Range("H136").Value = 1323543.4
$ws.Range("I136").Value = 1462771.9
$ws.Range("J136").Value = 873.25
$ws.Range("K136").Value = 4388315.699999999
$ws.Range("L136").Value = 2619.75
$ws.Range("M136").Value = -4385765.699999999
$ws.Range("N136").Value = -7719.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 818009.25
$ws.Range("I3").Value = 739.95654
$ws.Range("K3").Value = 739.95654
$ws.Range("M3").Value = -625.95654

$ws.Range("H94").Value = 879.55
$ws.Range("I94").Value = 792.73334
$ws.Range("K94").Value = 792.73334
$ws.Range("M94").Value = -341.73334

$ws.Range("H134").Value = 3409217.2
$ws.Range("I134").Value = 3528777.5
$ws.Range("J134").Value = 1750
$ws.Range("K134").Value = 10586332.5
$ws.Range("L134").Value = 5250
$ws.Range("M134").Value = -10583797.5
$ws.Range("N134").Value = -10320

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4734058
$ws.Range("I31").Value = 1261
$ws.Range("J31").Value = 10858854
$ws.Range("K31").Value = 1261
$ws.Range("L31").Value = 10858854
$ws.Range("M31").Value = -966
$ws.Range("N31").Value = -10859444

$ws.Range("H34").Value = 4734058
$ws.Range("I34").Value = 1261
$ws.Range("J34").Value = 10858854
$ws.Range("K34").Value = 1261
$ws.Range("L34").Value = 10858854
$ws.Range("M34").Value = -1059
$ws.Range("N34").Value = -10859258

$ws.Range("H50").Value = 10218.125
$ws.Range("J50").Value = 11680.4
$ws.Range("L50").Value = 11680.4
$ws.Range("N50").Value = -12930.4

$ws.Range("H58").Value = 1768.6234
$ws.Range("I58").Value = 844.9583
$ws.Range("K58").Value = 844.9583
$ws.Range("M58").Value = -641.9583

$ws.Range("H68").Value = 16411.445
$ws.Range("J68").Value = 19617.166
$ws.Range("L68").Value = 19617.166
$ws.Range("N68").Value = -21115.166

$ws.Range("H71").Value = 16411.445
$ws.Range("J71").Value = 19617.166
$ws.Range("L71").Value = 58851.49800000001
$ws.Range("N71").Value = -66339.49800000001

$ws.Range("H100").Value = 40779.5
$ws.Range("J100").Value = 40779.5
$ws.Range("L100").Value = 40779.5
$ws.Range("N100").Value = -42943.5

$ws.Range("H132").Value = 1605.36
$ws.Range("I132").Value = 1530.5227
$ws.Range("J132").Value = 2154.1667
$ws.Range("K132").Value = 4591.5681
$ws.Range("L132").Value = 6462.500100000001
$ws.Range("M132").Value = -2061.5681
$ws.Range("N132").Value = -11522.5001

$ws.Range("H136").Value = 1768.6234
$ws.Range("I136").Value = 844.9583
$ws.Range("K136").Value = 2534.8749
$ws.Range("M136").Value = 15.1251000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 453.51086
$ws.Range("I113").Value = 442.7931
$ws.Range("K113").Value = 1328.3793
$ws.Range("M113").Value = 841.6206999999999

$ws.Range("H131").Value = 3361.932
$ws.Range("J131").Value = 2503.7878
$ws.Range("L131").Value = 7511.3634
$ws.Range("N131").Value = -17591.3634

$ws.Range("H133").Value = 8620
$ws.Range("J133").Value = 9750
$ws.Range("L133").Value = 29250
$ws.Range("N133").Value = -39370

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1339.6522
$ws.Range("I82").Value = 1145.7142
$ws.Range("J82").Value = 1641.3334
$ws.Range("K82").Value = 1145.7142
$ws.Range("L82").Value = 1641.3334
$ws.Range("M82").Value = -784.7141999999999
$ws.Range("N82").Value = -2363.3334

$ws.Range("H85").Value = 1339.6522
$ws.Range("I85").Value = 1145.7142
$ws.Range("J85").Value = 1641.3334
$ws.Range("K85").Value = 1145.7142
$ws.Range("L85").Value = 1641.3334
$ws.Range("M85").Value = 102.2858000000001
$ws.Range("N85").Value = -4137.3334

$ws.Range("H100").Value = 3709.3635
$ws.Range("I100").Value = 1601.2
$ws.Range("J100").Value = 5466.1665
$ws.Range("K100").Value = 1601.2
$ws.Range("L100").Value = 5466.1665
$ws.Range("M100").Value = -1060.2
$ws.Range("N100").Value = -6548.1665

$ws.Range("H132").Value = 3826.561
$ws.Range("I132").Value = 3986.9487
$ws.Range("J132").Value = 699
$ws.Range("K132").Value = 11960.8461
$ws.Range("L132").Value = 2097
$ws.Range("M132").Value = -9430.846099999999
$ws.Range("N132").Value = -7157

$ws.Range("H136").Value = 1728
$ws.Range("I136").Value = 827.6
$ws.Range("J136").Value = 3728.889
$ws.Range("K136").Value = 2482.8
$ws.Range("L136").Value = 11186.667
$ws.Range("M136").Value = 67.19999999999982
$ws.Range("N136").Value = -16286.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7846786
$ws.Range("I132").Value = 8916667
$ws.Range("J132").Value = 990.5
$ws.Range("K132").Value = 26750001
$ws.Range("L132").Value = 2971.5
$ws.Range("M132").Value = -26747471
$ws.Range("N132").Value = -8031.5

$ws.Range("H136").Value = 2980873.5
$ws.Range("I136").Value = 5487.625
$ws.Range("J136").Value = 17857802
$ws.Range("K136").Value = 16462.875
$ws.Range("L136").Value = 53573406
$ws.Range("M136").Value = -13912.875
$ws.Range("N136").Value = -53578506
